$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Hunk 1: merge the two runs "-  " and "На странице " (row 8 / col 2 of the
# table) into a single run "-  На странице ". Scope the Find to that one
# table cell so the other 3 occurrences of "На странице " elsewhere in the
# document are left untouched.
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item(1)
$cell = $tbl.Cell(8, 2)
$cellRange = $cell.Range
$cellRange.Find.Execute("-  На странице ", $false, $false, $false, $false, $false, $true, 1, $false, "-  На странице ", 2) | Out-Null

# ---------------------------------------------------------------------------
# Hunk 2: remove the blank paragraph + the Russian "offtop" parenthetical
# remark at the very end of the document, keeping the paragraph that
# carries the (hidden) _GoBack bookmark, and flip that paragraph mark's
# language to en-US.
#
# Note: after the Find/Replace above, $d.Paragraphs (and Range objects
# fetched from it) can return stale/incorrect Start/End values for many
# paragraphs in this table-heavy document. Re-fetching the paragraph
# collection from $d.Content each time works around that.
# ---------------------------------------------------------------------------

$count = $d.Content.Paragraphs.Count
$blankPara = $d.Content.Paragraphs.Item($count - 1)  # empty paragraph right before the bookmark one

# 1) Delete the blank paragraph that precedes the bookmark paragraph.
$blankPara.Range.Delete() | Out-Null

# 2) Clear out the bookmark paragraph's runs/proofErr markup while
#    preserving the _GoBack bookmark. Text-range deletion alone leaves
#    orphaned <w:proofErr/> markers behind (they aren't addressable via
#    character ranges), so instead: split off a fresh trailing paragraph,
#    delete the old (bookmark-carrying) paragraph entirely -- which
#    cleanly drops its runs *and* its proofErr markers -- then re-plant
#    a _GoBack bookmark into the new paragraph.
$count = $d.Content.Paragraphs.Count
$target = $d.Content.Paragraphs.Item($count)
$target.Range.InsertParagraphAfter() | Out-Null

$count = $d.Content.Paragraphs.Count
$oldPara = $d.Content.Paragraphs.Item($count - 1)
$oldPara.Range.Delete() | Out-Null

# The engine mishandles Bookmarks.Add exactly at the absolute end of the
# document's content, so give it a throwaway character to anchor against,
# add the bookmark just before that character, then delete the character.
$count = $d.Content.Paragraphs.Count
$newLast = $d.Content.Paragraphs.Item($count)
$newLast.Range.InsertBefore("Z") | Out-Null

$count = $d.Content.Paragraphs.Count
$newLast = $d.Content.Paragraphs.Item($count)
$anchor = $d.Range($newLast.Range.Start, $newLast.Range.Start)
$d.Bookmarks.Add("_GoBack", $anchor) | Out-Null

$count = $d.Content.Paragraphs.Count
$newLast = $d.Content.Paragraphs.Item($count)
$placeholder = $d.Range($newLast.Range.End - 2, $newLast.Range.End - 1)
$placeholder.Delete() | Out-Null

# 3) Flip the (now run-less) paragraph mark's language to en-US, keeping
#    the existing bold formatting.
$count = $d.Content.Paragraphs.Count
$finalPara = $d.Content.Paragraphs.Item($count)
$finalPara.Range.LanguageID = "en-US"
